$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (data row): replace RAPHAELA's record with THEOMAR's record (new balance).
# Force the account number to be stored as text so the leading zeros survive
# (it would otherwise be auto-coerced to a number), then strip the
# quote-prefix formatting that this leaves behind so the cell matches the
# plain (unstyled) text cells used throughout the rest of the column.
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "004231509"
$ws.Range("A5").ClearFormats()
$ws.Range("B5").Value = "THEOMAR"
$ws.Range("C5").Value = 1326.97

# Row 46 held THEOMAR's old record (004231509 / THEOMAR / 413.97); it is
# removed entirely, shifting all subsequent rows up by one.
$ws.Rows(46).Delete()
